# Weekly Fruit/Vegetable price sheet: a new week's record is prepended to the
# "Femacal de La Calera - Ciboulette" price-history table (which starts at
# row 69, right after an earlier block of rows). Inserting a whole row there
# pushes the existing rows 69-318 down to 70-319 (Excel keeps every other
# column of those rows intact automatically), and we only need to populate
# the freshly inserted row 69 with this week's observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the first data row of the Ciboulette block; this
# shifts rows 69-318 down to 70-319 (and the sheet's UsedRange/dimension
# grows from R318 to R319) without touching any of their existing values.
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new week's record. Every
# column besides the date (Fecha) repeats the same constant values used
# throughout this market/category table.
$ws.Range("A69").Value = 3
$ws.Range("B69").Value = 'Femacal de La Calera'
$ws.Range("C69").Value = 'Coquimbo'
$ws.Range("D69").Value = 44707
$ws.Range("E69").Value = 5
$ws.Range("F69").Value = 100112039
$ws.Range("G69").Value = 'Ciboulette'
$ws.Range("H69").Value = 'Sin especificar'
$ws.Range("I69").Value = 'Primera'
$ws.Range("J69").Value = 180
$ws.Range("K69").Value = 1500
$ws.Range("L69").Value = 1500
$ws.Range("M69").Value = 1500
$ws.Range("N69").Value = '$/docena de atados'
$ws.Range("O69").Value = 'Provincia de Quillota'
$ws.Range("P69").Value = 500
$ws.Range("Q69").Value = 3
$ws.Range("R69").Value = 'Hortaliza'
